$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 361, pushing existing rows 361:483 down to 362:484
$ws.Rows("361:361").Insert()

# Populate the new row 361 with the new data point
$ws.Range("A361").Value = 5
$ws.Range("B361").Value = "Macroferia Regional de Talca"
$ws.Range("C361").Value = "Maule"
$ws.Range("D361").Value = 44988
$ws.Range("E361").Value = 7
$ws.Range("F361").Value = 100114013
$ws.Range("G361").Value = "Zanahoria"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 600
$ws.Range("K361").Value = 7000
$ws.Range("L361").Value = 7000
$ws.Range("M361").Value = 7000
$ws.Range("N361").Value = "$/saco 20 kilos"
$ws.Range("O361").Value = "Región de Ñuble"
$ws.Range("P361").Value = 350
$ws.Range("Q361").Value = 20
$ws.Range("R361").Value = "Hortaliza"
